# Eliminate colonne come da indicazioni
# Clears the placeholder "APPLICABILITA' = NO" / "RAZIONALE DI APPLICABILITA' = Tipo
# Documento non gestito" values from column J/K across the TestCases sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Rows whose J (APPLICABILITA') / K (RAZIONALE DI APPLICABILITA') columns hold the
# placeholder "NO" / "Tipo Documento non gestito" pair that needs to be removed.
$rows = @(8,9,10,11,12,13,14,15,16,17,18,19,20,21,23,24,25,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,105,106,107,108,109)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 10).ClearContents() # column J
    $ws.Cells.Item($r, 11).ClearContents() # column K
}

# Row 8 keeps a stray single space in column J (left behind by the clean-up).
$ws.Cells.Item(8, 10).Value2 = " "

# Restore the view state (scroll position / active selection) left behind after the edit.
$ws.Activate()
$ws.Range("E5").Select()
